# Apply the daily cryptos list refresh (coin prices / 1h volume %, and
# a handful of rows whose rank order shifted) as produced by the scraper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value reads as a plain number (e.g. "214.94")
# must be forced to Text format first, otherwise Excel auto-converts the
# assigned string into a numeric value and the literal text is lost -
# this mirrors the original scraped cells, which are all stored as text.
$textFormatCells = @(
    "D5",
    "D8",
    "D9",
    "D11",
    "D14",
    "D15",
    "D16",
    "D17",
    "D19",
    "D22",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D34",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D42",
    "D43",
    "D44",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Per-cell values for the refreshed row data.
$ws.Range("D2").Value = "27.022.21"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.676.26"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "214.94"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  +1.78%  "
$ws.Range("D9").Value = "21.29"
$ws.Range("E9").Value = "  +4.65%  "
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").Value = "0.0888"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "1.911.30"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "1.669.07"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "4.12"
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").Value = "0.535"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("D16").Value = "66.16"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "8.24"
$ws.Range("E17").Value = "  +3.06%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "27.010.42"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "235.29"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "0.0₃0738"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("D22").Value = "4.48"
$ws.Range("E22").Value = "  +1.38%  "
$ws.Range("D23").Value = "9.26"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("D25").Value = "146.30"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "7.25"
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").Value = "16.34"
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").Value = "0.113"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").Value = "1.536.88"
$ws.Range("E33").Value = "  +5.41%  "
$ws.Range("D34").Value = "3.17"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").Value = "0.590"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "0.913"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0174"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("D40").Value = "1.08"
$ws.Range("E40").Value = "  +9.08%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").Value = "67.53"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("D43").Value = "5.50"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("D44").Value = "2.26"
$ws.Range("E44").Value = "  -1.71%  "
$ws.Range("D45").Value = "1.819.26"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").Value = "90.51"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.54"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "8.06"
$ws.Range("E51").Value = "  +6.36%  "
